$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-04-23 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-04-24 Wednesday", 2)

$d.Content.Find.Execute("375×2=", $true, $false, $false, $false, $false, $true, 1, $false, "380×2=", 2)
$d.Content.Find.Execute("162×3=", $true, $false, $false, $false, $false, $true, 1, $false, "779×4=", 2)
$d.Content.Find.Execute("214×6=", $true, $false, $false, $false, $false, $true, 1, $false, "477×3=", 2)
$d.Content.Find.Execute("847×9=", $true, $false, $false, $false, $false, $true, 1, $false, "224×2=", 2)
$d.Content.Find.Execute("298×5=", $true, $false, $false, $false, $false, $true, 1, $false, "922×2=", 2)

$d.Content.Find.Execute("338×5=", $true, $false, $false, $false, $false, $true, 1, $false, "584×6=", 2)
$d.Content.Find.Execute("665×4=", $true, $false, $false, $false, $false, $true, 1, $false, "136×7=", 2)
$d.Content.Find.Execute("757×9=", $true, $false, $false, $false, $false, $true, 1, $false, "490×2=", 2)
$d.Content.Find.Execute("638×4=", $true, $false, $false, $false, $false, $true, 1, $false, "295×2=", 2)
$d.Content.Find.Execute("408×2=", $true, $false, $false, $false, $false, $true, 1, $false, "962×3=", 2)

$d.Content.Find.Execute("789×3=", $true, $false, $false, $false, $false, $true, 1, $false, "623×8=", 2)
$d.Content.Find.Execute("917×6=", $true, $false, $false, $false, $false, $true, 1, $false, "895×3=", 2)
$d.Content.Find.Execute("883×4=", $true, $false, $false, $false, $false, $true, 1, $false, "356×5=", 2)
$d.Content.Find.Execute("533×7=", $true, $false, $false, $false, $false, $true, 1, $false, "123×6=", 2)
$d.Content.Find.Execute("479×2=", $true, $false, $false, $false, $false, $true, 1, $false, "519×3=", 2)

$d.Content.Find.Execute("684×4=", $true, $false, $false, $false, $false, $true, 1, $false, "768×5=", 2)
$d.Content.Find.Execute("330×4=", $true, $false, $false, $false, $false, $true, 1, $false, "989×4=", 2)
$d.Content.Find.Execute("697×4=", $true, $false, $false, $false, $false, $true, 1, $false, "843×7=", 2)
$d.Content.Find.Execute("101×9=", $true, $false, $false, $false, $false, $true, 1, $false, "512×7=", 2)
$d.Content.Find.Execute("899×2=", $true, $false, $false, $false, $false, $true, 1, $false, "595×3=", 2)

$d.Content.Find.Execute("982×3=", $true, $false, $false, $false, $false, $true, 1, $false, "571×8=", 2)
$d.Content.Find.Execute("258×7=", $true, $false, $false, $false, $false, $true, 1, $false, "797×7=", 2)
$d.Content.Find.Execute("741×6=", $true, $false, $false, $false, $false, $true, 1, $false, "321×9=", 2)
$d.Content.Find.Execute("106×7=", $true, $false, $false, $false, $false, $true, 1, $false, "402×6=", 2)
$d.Content.Find.Execute("853×5=", $true, $false, $false, $false, $false, $true, 1, $false, "222×6=", 2)
